# Append the new 2025-10-04 row of portfolio data as row 50.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 50

# Column A holds the date as plain text (matches the existing A2:A49 cells,
# which are stored as text rather than real dates). Pre-format the cell as
# text so Excel doesn't silently convert the "2025-10-04" string into a date
# serial number, then clear the formatting back off again so the cell is
# left without any style override (same as its neighbors).
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-10-04"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = 54.45999908447266
$ws.Cells.Item($newRow, 3).Value = 716.0999755859375
$ws.Cells.Item($newRow, 4).Value = 328.4500122070312
